# DeleteProblemTestFactory tested on 2 lines of data
# Adds a second "problem" data row (row 3) to Sheet1, mirroring the layout
# of the existing row 2, and attaches a real hyperlink to the image-URL
# cell (H3), just like H2 would carry if it only had a single URL.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Write the new row of data (row 3) ------------------------------
$ws.Range("A3").Value2 = "2"
$ws.Range("B3").Value2 = "50.32"
$ws.Range("C3").Value2 = "30.32"
$ws.Range("D3").Value2 = "problemTitle2"
$ws.Range("E3").Value2 = "Сміттєзвалища"
$ws.Range("F3").Value2 = "problemDescription2"
$ws.Range("G3").Value2 = "problemSolution2"
$ws.Range("H3").Value2 = "http://i.imgur.com/1K6AdCH.jpg"
$ws.Range("I3").Value2 = "imageComment1"
$ws.Range("J3").Value2 = "admin@.com"
$ws.Range("K3").Value2 = "admin"

# --- 2. Turn the image-URL cell into a real hyperlink -------------------
$ws.Hyperlinks.Add($ws.Range("H3"), "http://i.imgur.com/1K6AdCH.jpg")

# --- 3. Re-apply the same cell format row 2 uses for its H/I columns ----
# (Hyperlinks.Add swaps in Excel's built-in "Hyperlink" style; row 2's
# H2/I2 cells instead use the sheet's own wrapped/top-aligned text style,
# so copy that formatting back onto H3/I3 to match.)
$ws.Range("I2").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("I3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
